# Applies the "Updated with version numbers" change:
#   - Insert a new column A named "Version" (shifts Code/Description/
#     Parent_Code/Parent_Description/Definition from A-E to B-F).
#   - Fill the new Version column with "1.0.0" for every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column at A, shifting existing columns A:E to B:F.
$ws.Columns.Item(1).Insert()

# Header for the new column.
$ws.Range("A1").Value = "Version"

# Determine the extent of the data (UsedRange now spans A1:F<lastRow>).
$lastRow = $ws.UsedRange.Rows.Count

# Fill the version number for every data row.
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 1).Value = "1.0.0"
}
